$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.01
$ws.Range("B7").Value = 5.910000000000001
$ws.Range("A10").Value = -21.663
$ws.Range("E10").Value = 16.364
$ws.Range("A12").Value = -21.606
$ws.Range("E14").Value = 17.037
$ws.Range("B15").Value = 5.072000000000001
$ws.Range("A18").Value = -22.002
$ws.Range("D18").Value = -8.638
$ws.Range("D19").Value = -7.912000000000001
$ws.Range("B20").Value = 7.633
$ws.Range("D27").Value = -8.183
$ws.Range("B29").Value = 5.645
$ws.Range("B30").Value = 5.694000000000001
$ws.Range("B31").Value = 6.276000000000001
$ws.Range("E32").Value = 16.701
$ws.Range("E35").Value = 16.194
$ws.Range("A37").Value = -19.92
$ws.Range("B40").Value = 9.327999999999999
$ws.Range("D42").Value = -8.257999999999999
$ws.Range("E43").Value = 17.139
$ws.Range("D44").Value = -8.125999999999999
$ws.Range("D47").Value = -7.815
$ws.Range("E49").Value = 16.359
$ws.Range("A55").Value = -21.841
$ws.Range("E56").Value = 16.276
$ws.Range("D58").Value = -8.370999999999999
$ws.Range("A68").Value = -21.50700000000001
$ws.Range("B68").Value = 5.881
$ws.Range("E69").Value = 17.448
$ws.Range("D73").Value = -8.047000000000001
$ws.Range("B76").Value = 6.341000000000001
$ws.Range("A77").Value = -20.637
$ws.Range("A78").Value = -20.126
$ws.Range("E81").Value = 16.782
$ws.Range("B87").Value = 4.462
$ws.Range("B88").Value = 5.051
$ws.Range("E92").Value = 17.834
$ws.Range("D95").Value = -7.567
$ws.Range("B96").Value = 6.468999999999999
$ws.Range("B98").Value = 5.502
$ws.Range("B101").Value = 7.87
$ws.Range("D101").Value = -8.010999999999999
$ws.Range("B102").Value = 7.739999999999999
